$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (image_number, filename, cell_1..cell_16)
$newRows = @(
    @(8, "008.jpg", 0,0,1,0,0,0,1,0,0,0,1,0,0,0,0,0),
    @(6, "006.jpg", 0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(3, "003.jpg", 0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(4, "004.jpg", 0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(-1, "pipeline_debug.jpg", 0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
}
